$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" value (row 8, column B)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: the two "Mapping" columns (AK = col 37, AL = col 38)
#    were reordered upstream, so both the header text and every data row
#    need their AK/AL contents swapped, and the column widths swap with them.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Header row (row 1): swap the two mapping header labels.
$ws.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR territorial division"
$ws.Range("AL1").Value = "Mapping: RIM Mapping"

# Data rows: swap AK<->AL per row (rows where both sides already match are
# left untouched since swapping them is a no-op).
$ws.Range("AK2").Value = "DivisionTerritoriale"
$ws.Range("AL2").Value = ""

$ws.Range("AK3").Value = ""
$ws.Range("AL3").Value = "n/a"

$ws.Range("AK6").Value = ""
$ws.Range("AL6").Value = "n/a"

$ws.Range("AK8").Value = ""
$ws.Range("AL8").Value = "N/A"

$ws.Range("AK9").Value = "type"
$ws.Range("AL9").Value = "N/A"

$ws.Range("AK11").Value = ""
$ws.Range("AL11").Value = "n/a"

$ws.Range("AK13").Value = ""
$ws.Range("AL13").Value = "N/A"

$ws.Range("AK14").Value = "code"
$ws.Range("AL14").Value = "N/A"

$ws.Range("AK15").Value = ""
$ws.Range("AL15").Value = "N/A"

$ws.Range("AK16").Value = ""
$ws.Range("AL16").Value = "N/A"

# Column widths: AK (37) and AL (38) swap their widths too.
$ws.Columns.Item(37).ColumnWidth = 72.98307291666667
$ws.Columns.Item(38).ColumnWidth = 24.147135416666668
